$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dimension 0 block (rows 5-7) ---
# Row 5: Average MST Weight - update n=2048,4096 values; drop n=8192,16834 columns
$ws.Range("L5").Value2 = 1.2023200000000001
$ws.Range("M5").Value2 = 1.20123
$ws.Range("N5:O5").ClearContents() | Out-Null

# Row 6: Max included edge - update n=2048 value; drop n=4096,8192,16834 columns
$ws.Range("L6").Value2 = 0.00658
$ws.Range("M6:O6").ClearContents() | Out-Null

# Row 7: (trials) - drop stray n=16834 trial count, add n=2048 trial count
$ws.Range("O7").ClearContents() | Out-Null
$ws.Range("L7").Value2 = 100

# --- Dimension 2 block (rows 11-13), extend out to n=65536 ---
# Row 11: Average MST Weight
$ws.Range("L11").Value2 = 29.677085999999999
$ws.Range("M11").Value2 = 41.798164
$ws.Range("N11").Value2 = 59.030811
$ws.Range("O11").Value2 = 84.327415000000002
$ws.Range("P11").Value2 = 117.4776
$ws.Range("Q11").Value2 = 166.05796799999999

# Row 12: Max included edge
$ws.Range("F12").Value2 = 0.44828499999999999
$ws.Range("L12").Value2 = 0.058646
$ws.Range("M12").Value2 = 0.033639
$ws.Range("N12").Value2 = 0.022433
$ws.Range("O12").Value2 = 0.015831
$ws.Range("P12").Value2 = 0.014917
$ws.Range("Q12").Value2 = 0.00867

# Row 13: (trials)
$ws.Range("M13").Value2 = 10
$ws.Range("N13").Value2 = 10
$ws.Range("O13").Value2 = 10
$ws.Range("P13").Value2 = 5
$ws.Range("Q13").Value2 = 5

# --- Remove leftover duplicate "NEW WITH ADJACENCY LIST" block and stray scratch rows ---
$ws.Rows("28:32").ClearContents() | Out-Null
$ws.Rows("36:39").ClearContents() | Out-Null

# --- Sheet view: scroll so column G is leftmost, select Q17 ---
$excel.Goto($ws.Range("G1"), $true) | Out-Null
$ws.Range("Q17").Select() | Out-Null

Write-Host "Done"
